# Update column F (dSF) values on Sheet1 to reflect repulled data / push all data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = -5
    3  = -2
    4  = -4
    5  = -2
    6  = -3
    7  = 9
    8  = 5
    9  = 5
    10 = 2
    12 = 3
    13 = 2
    14 = -3
    15 = -3
    16 = 1
    17 = -1
    18 = 1
    19 = -1
    20 = 3
    21 = -1
    22 = -4
    25 = -2
    26 = -5
    27 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
